$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 0.98702981104082932
$ws.Range("BP1").Value = 0.82623705142961223
$ws.Range("A2").Value = 0.83301819973320046
$ws.Range("Q2").Value = 0.94776298269716808
$ws.Range("V2").Value = 0.67928246348088916
$ws.Range("AW2").Value = 0.92955377476203327
$ws.Range("D3").Value = 0.66714454629032705
$ws.Range("BA4").Value = 0.96657455251813285
$ws.Range("B5").Value = 0.97739534477758516
$ws.Range("BK5").Value = 0.92144283652356196
$ws.Range("E6").Value = 0.94716171305632901
$ws.Range("H6").Value = 0.93485981683849273
$ws.Range("F7").Value = 0.89677664491672882
$ws.Range("I7").Value = 0.72249300366651181
$ws.Range("G8").Value = 0.71155010228448345
$ws.Range("I8").Value = 0.64964787247656175
$ws.Range("J9").Value = 0.92906798381968403
$ws.Range("AL9").Value = 0.78368070540622847
$ws.Range("H10").Value = 0.94478462470499069
$ws.Range("K10").Value = 0.91373327907225677
$ws.Range("AU10").Value = 0.69150846366307483
$ws.Range("I11").Value = 0.8342320277036408
$ws.Range("L11").Value = 0.81174551567362319
$ws.Range("C13").Value = 0.97210687429943543
$ws.Range("L13").Value = 0.77989137001795672
$ws.Range("L14").Value = 0.61975350650357219
$ws.Range("M14").Value = 0.54476965136685829
$ws.Range("AS14").Value = 0.58516142065028065
$ws.Range("M15").Value = 0.95802836424529136
$ws.Range("N15").Value = 0.92802981634508785
$ws.Range("P15").Value = 0.88811183255185888
$ws.Range("E16").Value = 0.87052521830084695
$ws.Range("W16").Value = 0.57603526640910263
$ws.Range("O17").Value = 0.77700702059295557
$ws.Range("T18").Value = 0.76603306474053057
$ws.Range("R19").Value = 0.85603729261816119
$ws.Range("K20").Value = 0.79731362970515063
$ws.Range("S20").Value = 0.95039062707563915
$ws.Range("S21").Value = 0.79112091435436627
$ws.Range("T21").Value = 0.97999741083878722
$ws.Range("BH21").Value = 0.71859879222378908
$ws.Range("T22").Value = 0.65666667581130611
$ws.Range("U22").Value = 0.66676795014831536
$ws.Range("S23").Value = 0.85238724763092599
$ws.Range("V24").Value = 0.96613419560497715
$ws.Range("W24").Value = 0.66292053889014613
$ws.Range("Y24").Value = 0.84889914156694557
$ws.Range("BN24").Value = 0.75966891379577128
$ws.Range("Z25").Value = 0.83936592414256384
$ws.Range("AA25").Value = 0.78320131695701867
$ws.Range("BA25").Value = 0.86830378431147237
$ws.Range("AA26").Value = 0.69097966779957209
$ws.Range("AC27").Value = 0.67252175901592559
$ws.Range("Z28").Value = 0.73134797423992781
$ws.Range("AA28").Value = 0.87174776349688876
$ws.Range("AD28").Value = 0.8168361260638014
$ws.Range("AZ28").Value = 0.91244101438822878
$ws.Range("AH29").Value = 0.96160949045900468
$ws.Range("D30").Value = 0.5845085615483645
$ws.Range("AG31").Value = 0.70018453773432854
$ws.Range("AD32").Value = 0.89213828311990484
$ws.Range("AE32").Value = 0.87857820767786454
$ws.Range("BA32").Value = 0.99360717272285881
$ws.Range("AC33").Value = 0.62913237187069337
$ws.Range("AF33").Value = 0.97562005097190529
$ws.Range("AI33").Value = 0.9451349825763411
$ws.Range("AG34").Value = 0.8896861066258821
$ws.Range("AI34").Value = 0.65045072713263674
$ws.Range("Q35").Value = 0.95426801858964194
$ws.Range("AK35").Value = 0.82972377508974593
$ws.Range("AH36").Value = 0.72651822786936648
$ws.Range("AI36").Value = 0.8420415505383374
$ws.Range("AJ37").Value = 0.88836854314742786
$ws.Range("AL37").Value = 0.6929569168713543
$ws.Range("BF37").Value = 0.82510142658196117
$ws.Range("K38").Value = 0.97894671700735547
$ws.Range("AJ38").Value = 0.83827440779751428
$ws.Range("AZ39").Value = 0.61858669769705932
$ws.Range("L40").Value = 0.59214638661559227
$ws.Range("AP40").Value = 0.99424032712788768
$ws.Range("AM41").Value = 0.71909215586619069
$ws.Range("AQ42").Value = 0.8812722974234497
$ws.Range("AR42").Value = 0.70328588205503739
$ws.Range("AI43").Value = 0.70023018000965265
$ws.Range("AO43").Value = 0.89464179881820871
$ws.Range("AR43").Value = 0.73055071379880543
$ws.Range("AV44").Value = 0.82756306218239439
$ws.Range("BJ44").Value = 0.81831459760107739
$ws.Range("AT45").Value = 0.86965931425251708
$ws.Range("BO45").Value = 0.89527681253007141
$ws.Range("AR46").Value = 0.93389533864361807
$ws.Range("AU46").Value = 0.71187627735389059
$ws.Range("AV46").Value = 0.95655796929883641
$ws.Range("AV47").Value = 0.71941675935656146
$ws.Range("BP48").Value = 0.96295165584063103
$ws.Range("AU49").Value = 0.85494943459614203
$ws.Range("AX49").Value = 0.80065520564080295
$ws.Range("AY49").Value = 0.83703930792555048
$ws.Range("AP50").Value = 0.7570087672885526
$ws.Range("AV50").Value = 0.98936897021034742
$ws.Range("AY50").Value = 0.97910130482975222
$ws.Range("AO51").Value = 0.84054981642780291
$ws.Range("AY52").Value = 0.68142271431292223
$ws.Range("AY53").Value = 0.9354947295749555
$ws.Range("AZ53").Value = 0.95297108712631529
$ws.Range("AZ54").Value = 0.76550312846395319
$ws.Range("BC54").Value = 0.6239340317689841
$ws.Range("BD54").Value = 0.9232758607249848
$ws.Range("BD55").Value = 0.75520023351415078
$ws.Range("BE55").Value = 0.94669916851226144
$ws.Range("BE56").Value = 0.95831493640043885
$ws.Range("BF56").Value = 0.7504443008206132
$ws.Range("AE57").Value = 0.83364886124462734
$ws.Range("BF57").Value = 0.83607376631026153
$ws.Range("AI58").Value = 0.88754641788814226
$ws.Range("BF59").Value = 0.87469615267175282
$ws.Range("BG60").Value = 0.94566145180306971
$ws.Range("BI60").Value = 0.90815882463315778
$ws.Range("BG61").Value = 0.98294321235426607
$ws.Range("BK61").Value = 0.86208532116581327
$ws.Range("BH62").Value = 0.89495581955735204
$ws.Range("BI62").Value = 0.93745613843225639
$ws.Range("BL62").Value = 0.71335936538759248
$ws.Range("B63").Value = 0.8510461504852701
$ws.Range("BK64").Value = 0.96417058223217844
$ws.Range("BL65").Value = 0.87507090810340138
$ws.Range("BO65").Value = 0.78236371404939131
$ws.Range("AR66").Value = 0.98459241290921229
$ws.Range("BA66").Value = 0.50155228117638551
$ws.Range("BL66").Value = 0.80189063397282223
$ws.Range("BM66").Value = 0.77118438891128149
$ws.Range("A67").Value = 0.9797204760663254
$ws.Range("BN67").Value = 0.93075959429579291
$ws.Range("BP67").Value = 0.83878487066512331
$ws.Range("D68").Value = 0.81171471832210829
$ws.Range("BN68").Value = 0.89628739025599569
